$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '47.540.03'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.25%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.494.96'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '322.08'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '108.92'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.55%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.26%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '39.30'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.38%  '
$ws.Range('E12').Value = '  +0.55%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '18.68'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('E14').Value = '  +0.48%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '2.883.16'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.17%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.491.99'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('E17').Value = '  +0.37%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '47.427.30'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.21%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.38'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +4.99%  '
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('E21').Value = '  +0.32%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +14.88%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '70.64'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.13%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '246.79'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.74%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.56'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('E28').Value = '  -0.45%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.139'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +4.02%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '34.66'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('E31').Value = '  -5.59%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '49.94'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.97%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '20.27'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('E34').Value = '  -0.41%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.0789'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.69%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.05%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.73'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.11%  '
$ws.Range('E38').Value = '  -0.08%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.96'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('E40').Value = '  +0.15%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '22.54'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +6.60%  '
$ws.Range('E42').Value = '  -2.16%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '118.95'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('E44').Value = '  -0.25%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.001.00'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.78%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.03'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('E49').Value = '  -1.99%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.45%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '56.64'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.23%  '
